$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "first observed import" data rows added for United Arab Emirates (row 22),
# Jordan (row 23) and Cyprus (row 25).
# Fill the date column (C) first, then the iso_code column (B), then the
# obs_imp_flag column (D), matching the order the data was originally entered.

$ws.Range("C25").Value = "2024-12-15"
$ws.Range("C23").Value = "2024-09-01"
$ws.Range("C22").Value = "2023-01-11"

$ws.Range("B22").Value = "ARE"
$ws.Range("B23").Value = "JOR"
$ws.Range("B25").Value = "CYP"

$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("D25").Value = 1

$ws.Range("F14").Select()
